$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 550
$ws.Range("J39").Value = 1000
$ws.Range("L39").Value = 3000
$ws.Range("N39").Value = -3592
$ws.Range("H64").Value = 15153496
$ws.Range("J64").Value = 35718424
$ws.Range("L64").Value = 35718424
$ws.Range("N64").Value = -35718920
$ws.Range("H67").Value = 15153496
$ws.Range("J67").Value = 35718424
$ws.Range("L67").Value = 35718424
$ws.Range("N67").Value = -35720140
$ws.Range("H107").Value = 186.5
$ws.Range("I107").Value = 135
$ws.Range("K107").Value = 135
$ws.Range("M107").Value = 1785
$ws.Range("H116").Value = 38225040
$ws.Range("I116").Value = 38028620
$ws.Range("J116").Value = 38466790
$ws.Range("K116").Value = 38028620
$ws.Range("L116").Value = 38466790
$ws.Range("M116").Value = -38025178
$ws.Range("N116").Value = -38473674
$ws.Range("H135").Value = 3937.1333
$ws.Range("I135").Value = 1518.7
$ws.Range("K135").Value = 13668.3
$ws.Range("M135").Value = -11133.3
$ws.Range("H137").Value = 12721621
$ws.Range("I137").Value = 1668733.4
$ws.Range("K137").Value = 5006200.199999999
$ws.Range("M137").Value = -5003650.199999999
$ws.Range("H138").Value = 5176.6
$ws.Range("J138").Value = 5759.1704
$ws.Range("L138").Value = 17277.5112
$ws.Range("N138").Value = -27557.5112
$ws.Range("H140").Value = 63132.25
$ws.Range("J140").Value = 62049.855
$ws.Range("L140").Value = 62049.855
$ws.Range("N140").Value = -72409.85500000001
$ws.Range("H141").Value = 9116.909
$ws.Range("I141").Value = 9638.6
$ws.Range("K141").Value = 28915.8
$ws.Range("M141").Value = -23735.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4236
$ws.Range("I32").Value = 2231.6667
$ws.Range("J32").Value = 15059.4
$ws.Range("K32").Value = 2231.6667
$ws.Range("L32").Value = 15059.4
$ws.Range("M32").Value = -1944.6667
$ws.Range("N32").Value = -15633.4
$ws.Range("H61").Value = 9483.9
$ws.Range("I61").Value = 11104.875
$ws.Range("K61").Value = 11104.875
$ws.Range("M61").Value = -10892.875
$ws.Range("H74").Value = 35716644
$ws.Range("I74").Value = 62501500
$ws.Range("K74").Value = 62501500
$ws.Range("M74").Value = -62500626
$ws.Range("H77").Value = 35716644
$ws.Range("I77").Value = 62501500
$ws.Range("K77").Value = 312507500
$ws.Range("M77").Value = -312503132
$ws.Range("H132").Value = 29939.479
$ws.Range("I132").Value = 48984.25
$ws.Range("K132").Value = 146952.75
$ws.Range("M132").Value = -144422.75
$ws.Range("H136").Value = 9483.9
$ws.Range("I136").Value = 11104.875
$ws.Range("K136").Value = 33314.625
$ws.Range("M136").Value = -30764.625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1398.5714
$ws.Range("J86").Value = 1356.5714
$ws.Range("L86").Value = 1356.5714
$ws.Range("N86").Value = -3602.5714
$ws.Range("H89").Value = 1398.5714
$ws.Range("J89").Value = 1356.5714
$ws.Range("L89").Value = 6782.857
$ws.Range("N89").Value = -18014.857
$ws.Range("H107").Value = 1720
$ws.Range("I107").Value = 1858.1818
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 1858.1818
$ws.Range("L107").Value = 200
$ws.Range("M107").Value = 61.81819999999993
$ws.Range("N107").Value = -4040
$ws.Range("H134").Value = 3614.4119
$ws.Range("I134").Value = 1147
$ws.Range("K134").Value = 3441
$ws.Range("M134").Value = -906

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2137.3333
$ws.Range("I16").Value = 2755.5
$ws.Range("J16").Value = 901
$ws.Range("K16").Value = 2755.5
$ws.Range("L16").Value = 901
$ws.Range("M16").Value = -2468.5
$ws.Range("N16").Value = -1475
$ws.Range("H31").Value = 6116.0757
$ws.Range("I31").Value = 2301.8333
$ws.Range("K31").Value = 2301.8333
$ws.Range("M31").Value = -2006.8333
$ws.Range("H34").Value = 6116.0757
$ws.Range("I34").Value = 2301.8333
$ws.Range("K34").Value = 2301.8333
$ws.Range("M34").Value = -2099.8333
$ws.Range("H105").Value = 1500
$ws.Range("H113").Value = 2137.3333
$ws.Range("I113").Value = 2755.5
$ws.Range("J113").Value = 901
$ws.Range("K113").Value = 2755.5
$ws.Range("L113").Value = 901
$ws.Range("M113").Value = -585.5
$ws.Range("N113").Value = -5241

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 226.33333
$ws.Range("I2").Value = 162
$ws.Range("J2").Value = 284.81818
$ws.Range("K2").Value = 972
$ws.Range("L2").Value = 1708.90908
$ws.Range("M2").Value = -859
$ws.Range("N2").Value = -1934.90908
$ws.Range("H4").Value = 3171376.5
$ws.Range("I4").Value = 1040431.9
$ws.Range("K4").Value = 3121295.7
$ws.Range("M4").Value = -3121183.7
$ws.Range("H5").Value = 243936.94
$ws.Range("I5").Value = 525.75
$ws.Range("K5").Value = 1577.25
$ws.Range("M5").Value = -1465.25
$ws.Range("H37").Value = 166762480
$ws.Range("J37").Value = 166762480
$ws.Range("L37").Value = 500287440
$ws.Range("N37").Value = -500287664
$ws.Range("H68").Value = 973744
$ws.Range("I68").Value = 1566.6666
$ws.Range("K68").Value = 4699.9998
$ws.Range("M68").Value = -3888.9998
$ws.Range("H71").Value = 973744
$ws.Range("I71").Value = 1566.6666
$ws.Range("K71").Value = 14099.9994
$ws.Range("M71").Value = -10043.9994
$ws.Range("H129").Value = 511
$ws.Range("I129").Value = 511
$ws.Range("K129").Value = 1533
$ws.Range("M129").Value = 3467
$ws.Range("H135").Value = 243936.94
$ws.Range("I135").Value = 525.75
$ws.Range("K135").Value = 4731.75
$ws.Range("M135").Value = -2196.75
$ws.Range("H140").Value = 3947.4688
$ws.Range("I140").Value = 2688.1304
$ws.Range("K140").Value = 8064.3912
$ws.Range("M140").Value = -2884.3912

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 607279.8
$ws.Range("I14").Value = 1500950
$ws.Range("J14").Value = 11499.667
$ws.Range("K14").Value = 1500950
$ws.Range("L14").Value = 11499.667
$ws.Range("M14").Value = -1500782
$ws.Range("N14").Value = -11835.667
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5302
$ws.Range("H132").Value = 99434.86
$ws.Range("I132").Value = 136916.53
$ws.Range("K132").Value = 410749.59
$ws.Range("M132").Value = -408219.59

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4819.8
$ws.Range("I7").Value = 3039.2
$ws.Range("J7").Value = 6600.4
$ws.Range("K7").Value = 3039.2
$ws.Range("L7").Value = 6600.4
$ws.Range("M7").Value = -2927.2
$ws.Range("N7").Value = -6824.4
$ws.Range("H18").Value = 14999
$ws.Range("J18").Value = 14999
$ws.Range("L18").Value = 14999
$ws.Range("N18").Value = -15343
$ws.Range("H22").Value = 1188.3889
$ws.Range("I22").Value = 999.8889
$ws.Range("J22").Value = 1376.8889
$ws.Range("K22").Value = 999.8889
$ws.Range("L22").Value = 1376.8889
$ws.Range("M22").Value = -704.8889
$ws.Range("N22").Value = -1966.8889
$ws.Range("H27").Value = 1188.3889
$ws.Range("I27").Value = 999.8889
$ws.Range("J27").Value = 1376.8889
$ws.Range("K27").Value = 999.8889
$ws.Range("L27").Value = 1376.8889
$ws.Range("M27").Value = -892.8889
$ws.Range("N27").Value = -1590.8889
$ws.Range("H46").Value = 5796.4688
$ws.Range("I46").Value = 3666.6667
$ws.Range("J46").Value = 6016.793
$ws.Range("K46").Value = 3666.6667
$ws.Range("L46").Value = 6016.793
$ws.Range("M46").Value = -3478.6667
$ws.Range("N46").Value = -6392.793
$ws.Range("H126").Value = 4819.8
$ws.Range("I126").Value = 3039.2
$ws.Range("J126").Value = 6600.4
$ws.Range("K126").Value = 9117.599999999999
$ws.Range("L126").Value = 19801.2
$ws.Range("M126").Value = -6647.599999999999
$ws.Range("N126").Value = -24741.2
$ws.Range("H132").Value = 5036.5713
$ws.Range("I132").Value = 5704.8
$ws.Range("K132").Value = 17114.4
$ws.Range("M132").Value = -14584.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6499.8
$ws.Range("I2").Value = 6730.5386
$ws.Range("K2").Value = 6730.5386
$ws.Range("M2").Value = -6618.5386
$ws.Range("H29").Value = 8999
$ws.Range("I29").Value = 8999
$ws.Range("K29").Value = 8999
$ws.Range("M29").Value = -8709
$ws.Range("H81").Value = 2989745.8
$ws.Range("I81").Value = 3476872
$ws.Range("K81").Value = 6953744
$ws.Range("M81").Value = -6952683
$ws.Range("H84").Value = 2989745.8
$ws.Range("I84").Value = 3476872
$ws.Range("K84").Value = 34768720
$ws.Range("M84").Value = -34763416
$ws.Range("H100").Value = 909701
$ws.Range("I100").Value = 1438412.6
$ws.Range("J100").Value = 3338.1428
$ws.Range("K100").Value = 2876825.2
$ws.Range("L100").Value = 6676.2856
$ws.Range("M100").Value = -2876284.2
$ws.Range("N100").Value = -7758.2856
$ws.Range("H126").Value = 2439.4736
$ws.Range("I126").Value = 2146.5
$ws.Range("K126").Value = 6439.5
$ws.Range("M126").Value = -3969.5
$ws.Range("H132").Value = 32682056
$ws.Range("I132").Value = 5556876
$ws.Range("K132").Value = 16670628
$ws.Range("M132").Value = -16668098
